$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark currently at the end of the title paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Append the new "Info" section (empty paragraph, bold heading, body
#    paragraph with the demo description) right before the final section
#    break, and re-create the _GoBack bookmark at the very end of the
#    inserted content.
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>' +
       '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
         '<w:pPr><w:rPr><w:b/><w:sz w:val="24"/></w:rPr></w:pPr>' +
         '<w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t>Info</w:t></w:r>' +
       '</w:p>' +
       '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
         '<w:pPr><w:ind w:firstLine="720"/></w:pPr>' +
         '<w:r><w:t>Pentru un stadiu intermediar, am realizat un scurt demo pentru aplicatie care deocamdata contine doar elemente pe care le vom folosi mai apoi in crearea aplicatiilor (ex. Pagina de log-in, Pagina cu harta, Pagina pentru client). Acest demo se afla in repository</w:t></w:r>' +
         '<w:r><w:t xml:space="preserve">: </w:t></w:r>' +
         '<w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="F6F8FA"/></w:rPr><w:t>Demo - app.mp4</w:t></w:r>' +
         '<w:r><w:t xml:space="preserve"> .</w:t></w:r>' +
         '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
       '</w:p>'

$endRange.InsertXML($xml)
